$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - replace with new NDC / medicine record (header row stays untouched)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "6275651818"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "Sun Pharmaceutical Industries, Inc."
$ws.Range("C2").Value = "CARBIDOPA AND LEVODOPA"
$ws.Range("D2").Value = "25 mg/1"
$ws.Range("E2").Value = "HAD1849A"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "04/24/30"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = 1

# Row 3 - replace with new NDC / medicine record
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2315574603"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "Avet Pharmaceuticals Inc."
$ws.Range("C3").Value = "Rasagiline mesylate"
$ws.Range("D3").Value = ".5 mg/1"
$ws.Range("E3").Value = "RCY01AD6"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "12/23/31"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = 1

# Row 4 - touch the row so an (empty) row entry is added beneath the data,
# extending the sheet's used range without writing any cell content to it.
$r = $ws.Rows.Item(4)
$r.OutlineLevel = 0
